$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test case 13 completed: two new rows of effort-log data, plus the
# matching "Task" description strings, then move the selection.
$ws.Range("A17").Value = 41445
$ws.Range("B17").Value = 0.5
$ws.Range("D17").Value = "Implementation tc13_eventStates"

$ws.Range("A18").Value = 41446
$ws.Range("C18").Value = 2.25
$ws.Range("D18").Value = "Successful completion of tc13"

$ws.Range("E11").Select()
